$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('O2').NumberFormat = "@"
$ws.Range('O2').Value = '2022-08-18 20:59:13'
$ws.Range('A3').NumberFormat = "@"
$ws.Range('A3').Value = '6075745012'
$ws.Range('B3').NumberFormat = "@"
$ws.Range('B3').Value = 'Avela Strumpfhose Madame Natural  11 - 12'
$ws.Range('C3').NumberFormat = "@"
$ws.Range('C3').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-madame-natural-11-12/p/6075745012'
$ws.Range('D3').Value = ""
$ws.Range('H3').NumberFormat = "@"
$ws.Range('H3').Value = '5.95'
$ws.Range('I3').Value = ""
$ws.Range('J3').Value = ""
$ws.Range('K3').Value = ""
$ws.Range('L3').Value = ""
$ws.Range('M3').NumberFormat = "@"
$ws.Range('M3').Value = '[''haushalt-tier'', ''bekleidung'', ''socken-unterwaesche'', ''struempfe'']'
$ws.Range('N3').NumberFormat = "@"
$ws.Range('N3').Value = 'Avela Strumpfhose Madame Natural  11 - 12 5.95 Schweizer Franken'
$ws.Range('O3').NumberFormat = "@"
$ws.Range('O3').Value = '2022-08-18 20:59:13'
$ws.Range('O4').NumberFormat = "@"
$ws.Range('O4').Value = '2022-08-18 20:59:13'
$ws.Range('O5').NumberFormat = "@"
$ws.Range('O5').Value = '2022-08-18 20:59:13'
$ws.Range('O6').NumberFormat = "@"
$ws.Range('O6').Value = '2022-08-18 20:59:13'
$ws.Range('O7').NumberFormat = "@"
$ws.Range('O7').Value = '2022-08-18 20:59:13'
$ws.Range('O8').NumberFormat = "@"
$ws.Range('O8').Value = '2022-08-18 20:59:13'
$ws.Range('O9').NumberFormat = "@"
$ws.Range('O9').Value = '2022-08-18 20:59:13'
$ws.Range('O10').NumberFormat = "@"
$ws.Range('O10').Value = '2022-08-18 20:59:13'
$ws.Range('O11').NumberFormat = "@"
$ws.Range('O11').Value = '2022-08-18 20:59:13'
$ws.Range('O12').NumberFormat = "@"
$ws.Range('O12').Value = '2022-08-18 20:59:13'
$ws.Range('O13').NumberFormat = "@"
$ws.Range('O13').Value = '2022-08-18 20:59:13'
$ws.Range('O14').NumberFormat = "@"
$ws.Range('O14').Value = '2022-08-18 20:59:13'
$ws.Range('A15').NumberFormat = "@"
$ws.Range('A15').Value = '6548194'
$ws.Range('B15').NumberFormat = "@"
$ws.Range('B15').Value = 'Selenacare Menstruationsunterwäsche S'
$ws.Range('C15').NumberFormat = "@"
$ws.Range('C15').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/unterwaesche/selenacare-menstruationsunterwaesche-s/p/6548194'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '1ST'
$ws.Range('E15').Value = ""
$ws.Range('F15').Value = 0
$ws.Range('G15').NumberFormat = "@"
$ws.Range('G15').Value = 'Selenacare'
$ws.Range('H15').NumberFormat = "@"
$ws.Range('H15').Value = '12.25'
$ws.Range('I15').NumberFormat = "@"
$ws.Range('I15').Value = '12.25/1ST'
$ws.Range('J15').NumberFormat = "@"
$ws.Range('J15').Value = 'Preis pro 1 Stück'
$ws.Range('K15').NumberFormat = "@"
$ws.Range('K15').Value = '12.25'
$ws.Range('L15').NumberFormat = "@"
$ws.Range('L15').Value = '1ST'
$ws.Range('M15').NumberFormat = "@"
$ws.Range('M15').Value = '[''haushalt-tier'', ''bekleidung'', ''socken-unterwaesche'', ''unterwaesche'']'
$ws.Range('N15').NumberFormat = "@"
$ws.Range('N15').Value = 'Selenacare Menstruationsunterwäsche S 50% Aktion 12.25 Schweizer Franken statt 24.50 Schweizer Franken'
$ws.Range('O15').NumberFormat = "@"
$ws.Range('O15').Value = '2022-08-18 20:59:13'
$ws.Range('A16').NumberFormat = "@"
$ws.Range('A16').Value = '4091029003'
$ws.Range('B16').NumberFormat = "@"
$ws.Range('B16').Value = 'Avela Socken Cotton Euro 35 - 38'
$ws.Range('C16').NumberFormat = "@"
$ws.Range('C16').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/socken/avela-socken-cotton-euro-35-38/p/4091029003'
$ws.Range('E16').Value = 1
$ws.Range('F16').Value = 5
$ws.Range('N16').NumberFormat = "@"
$ws.Range('N16').Value = 'Avela Socken Cotton Euro 35 - 38 5.50 Schweizer Franken'
$ws.Range('O16').NumberFormat = "@"
$ws.Range('O16').Value = '2022-08-18 20:59:13'
$ws.Range('A17').NumberFormat = "@"
$ws.Range('A17').Value = '4091029004'
$ws.Range('B17').NumberFormat = "@"
$ws.Range('B17').Value = 'Avela Socken Cotton Euro 39 - 41'
$ws.Range('C17').NumberFormat = "@"
$ws.Range('C17').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/socken/avela-socken-cotton-euro-39-41/p/4091029004'
$ws.Range('E17').Value = ""
$ws.Range('F17').Value = 0
$ws.Range('N17').NumberFormat = "@"
$ws.Range('N17').Value = 'Avela Socken Cotton Euro 39 - 41 5.50 Schweizer Franken'
$ws.Range('O17').NumberFormat = "@"
$ws.Range('O17').Value = '2022-08-18 20:59:13'
$ws.Range('A18').NumberFormat = "@"
$ws.Range('A18').Value = '4091029007'
$ws.Range('B18').NumberFormat = "@"
$ws.Range('B18').Value = 'Avela Socken Cotton Noir 35 - 38'
$ws.Range('C18').NumberFormat = "@"
$ws.Range('C18').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/socken/avela-socken-cotton-noir-35-38/p/4091029007'
$ws.Range('E18').Value = 1
$ws.Range('F18').Value = 2
$ws.Range('N18').NumberFormat = "@"
$ws.Range('N18').Value = 'Avela Socken Cotton Noir 35 - 38 5.50 Schweizer Franken'
$ws.Range('O18').NumberFormat = "@"
$ws.Range('O18').Value = '2022-08-18 20:59:13'
$ws.Range('A19').NumberFormat = "@"
$ws.Range('A19').Value = '4091029008'
$ws.Range('B19').NumberFormat = "@"
$ws.Range('B19').Value = 'Avela Socken Cotton Noir 39 - 41'
$ws.Range('C19').NumberFormat = "@"
$ws.Range('C19').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/socken/avela-socken-cotton-noir-39-41/p/4091029008'
$ws.Range('D19').Value = ""
$ws.Range('E19').Value = 2
$ws.Range('F19').Value = 5
$ws.Range('G19').NumberFormat = "@"
$ws.Range('G19').Value = 'Coop'
$ws.Range('H19').NumberFormat = "@"
$ws.Range('H19').Value = '5.50'
$ws.Range('I19').Value = ""
$ws.Range('J19').Value = ""
$ws.Range('K19').Value = ""
$ws.Range('L19').Value = ""
$ws.Range('N19').NumberFormat = "@"
$ws.Range('N19').Value = 'Avela Socken Cotton Noir 39 - 41 5.50 Schweizer Franken'
$ws.Range('O19').NumberFormat = "@"
$ws.Range('O19').Value = '2022-08-18 20:59:13'
$ws.Range('A20').NumberFormat = "@"
$ws.Range('A20').Value = '6077158004'
$ws.Range('B20').NumberFormat = "@"
$ws.Range('B20').Value = 'Avela Söckchen Ideal Hasel One Size'
$ws.Range('C20').NumberFormat = "@"
$ws.Range('C20').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/socken/avela-soeckchen-ideal-hasel-one-size/p/6077158004'
$ws.Range('N20').NumberFormat = "@"
$ws.Range('N20').Value = 'Avela Söckchen Ideal Hasel One Size 4.95 Schweizer Franken'
$ws.Range('O20').NumberFormat = "@"
$ws.Range('O20').Value = '2022-08-18 20:59:13'
$ws.Range('A21').NumberFormat = "@"
$ws.Range('A21').Value = '6077158006'
$ws.Range('B21').NumberFormat = "@"
$ws.Range('B21').Value = 'Avela Söckchen Ideal Noir One Size'
$ws.Range('C21').NumberFormat = "@"
$ws.Range('C21').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/socken/avela-soeckchen-ideal-noir-one-size/p/6077158006'
$ws.Range('H21').NumberFormat = "@"
$ws.Range('H21').Value = '4.95'
$ws.Range('I21').NumberFormat = "@"
$ws.Range('I21').Value = '2.48/1ST'
$ws.Range('K21').NumberFormat = "@"
$ws.Range('K21').Value = '2.48'
$ws.Range('N21').NumberFormat = "@"
$ws.Range('N21').Value = 'Avela Söckchen Ideal Noir One Size 4.95 Schweizer Franken'
$ws.Range('O21').NumberFormat = "@"
$ws.Range('O21').Value = '2022-08-18 20:59:13'
$ws.Range('A22').NumberFormat = "@"
$ws.Range('A22').Value = '6077154003'
$ws.Range('B22').NumberFormat = "@"
$ws.Range('B22').Value = 'Avela Söckchen Pure Hasel One Size'
$ws.Range('C22').NumberFormat = "@"
$ws.Range('C22').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/socken/avela-soeckchen-pure-hasel-one-size/p/6077154003'
$ws.Range('N22').NumberFormat = "@"
$ws.Range('N22').Value = 'Avela Söckchen Pure Hasel One Size 2.95 Schweizer Franken'
$ws.Range('O22').NumberFormat = "@"
$ws.Range('O22').Value = '2022-08-18 20:59:13'
$ws.Range('A23').NumberFormat = "@"
$ws.Range('A23').Value = '6077154004'
$ws.Range('B23').NumberFormat = "@"
$ws.Range('B23').Value = 'Avela Söckchen Pure Natural One Size'
$ws.Range('C23').NumberFormat = "@"
$ws.Range('C23').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/socken/avela-soeckchen-pure-natural-one-size/p/6077154004'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2ST'
$ws.Range('E23').Value = ""
$ws.Range('F23').Value = 0
$ws.Range('H23').NumberFormat = "@"
$ws.Range('H23').Value = '2.95'
$ws.Range('I23').NumberFormat = "@"
$ws.Range('I23').Value = '1.48/1ST'
$ws.Range('J23').NumberFormat = "@"
$ws.Range('J23').Value = 'Preis pro 1 Stück'
$ws.Range('K23').NumberFormat = "@"
$ws.Range('K23').Value = '1.48'
$ws.Range('L23').NumberFormat = "@"
$ws.Range('L23').Value = '1ST'
$ws.Range('M23').NumberFormat = "@"
$ws.Range('M23').Value = '[''haushalt-tier'', ''bekleidung'', ''socken-unterwaesche'', ''socken'']'
$ws.Range('N23').NumberFormat = "@"
$ws.Range('N23').Value = 'Avela Söckchen Pure Natural One Size 2.95 Schweizer Franken'
$ws.Range('O23').NumberFormat = "@"
$ws.Range('O23').Value = '2022-08-18 20:59:13'
$ws.Range('A24').NumberFormat = "@"
$ws.Range('A24').Value = '6075681007'
$ws.Range('B24').NumberFormat = "@"
$ws.Range('B24').Value = 'Avela Strumpfhose Belform Liberty Hasel  10 - 10.5'
$ws.Range('C24').NumberFormat = "@"
$ws.Range('C24').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-belform-liberty-hasel-10-105/p/6075681007'
$ws.Range('N24').NumberFormat = "@"
$ws.Range('N24').Value = 'Avela Strumpfhose Belform Liberty Hasel  10 - 10.5 6.95 Schweizer Franken'
$ws.Range('O24').NumberFormat = "@"
$ws.Range('O24').Value = '2022-08-18 20:59:13'
$ws.Range('A25').NumberFormat = "@"
$ws.Range('A25').Value = '6075681008'
$ws.Range('B25').NumberFormat = "@"
$ws.Range('B25').Value = 'Avela Strumpfhose Belform Liberty Hasel  11 - 12'
$ws.Range('C25').NumberFormat = "@"
$ws.Range('C25').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-belform-liberty-hasel-11-12/p/6075681008'
$ws.Range('E25').Value = 1
$ws.Range('F25').Value = 5
$ws.Range('N25').NumberFormat = "@"
$ws.Range('N25').Value = 'Avela Strumpfhose Belform Liberty Hasel  11 - 12 6.95 Schweizer Franken'
$ws.Range('O25').NumberFormat = "@"
$ws.Range('O25').Value = '2022-08-18 20:59:13'
$ws.Range('A26').NumberFormat = "@"
$ws.Range('A26').Value = '6075681005'
$ws.Range('B26').NumberFormat = "@"
$ws.Range('B26').Value = 'Avela Strumpfhose Belform Liberty Hasel  8.5 - 9'
$ws.Range('C26').NumberFormat = "@"
$ws.Range('C26').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-belform-liberty-hasel-85-9/p/6075681005'
$ws.Range('N26').NumberFormat = "@"
$ws.Range('N26').Value = 'Avela Strumpfhose Belform Liberty Hasel  8.5 - 9 6.95 Schweizer Franken'
$ws.Range('O26').NumberFormat = "@"
$ws.Range('O26').Value = '2022-08-18 20:59:13'
$ws.Range('A27').NumberFormat = "@"
$ws.Range('A27').Value = '6075681006'
$ws.Range('B27').NumberFormat = "@"
$ws.Range('B27').Value = 'Avela Strumpfhose Belform Liberty Hasel  9.5'
$ws.Range('C27').NumberFormat = "@"
$ws.Range('C27').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-belform-liberty-hasel-95/p/6075681006'
$ws.Range('E27').Value = ""
$ws.Range('F27').Value = 0
$ws.Range('N27').NumberFormat = "@"
$ws.Range('N27').Value = 'Avela Strumpfhose Belform Liberty Hasel  9.5 6.95 Schweizer Franken'
$ws.Range('O27').NumberFormat = "@"
$ws.Range('O27').Value = '2022-08-18 20:59:13'
$ws.Range('A28').NumberFormat = "@"
$ws.Range('A28').Value = '6075681011'
$ws.Range('B28').NumberFormat = "@"
$ws.Range('B28').Value = 'Avela Strumpfhose Belform Liberty Natural  10 - 10.5'
$ws.Range('C28').NumberFormat = "@"
$ws.Range('C28').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-belform-liberty-natural-10-105/p/6075681011'
$ws.Range('E28').Value = 1
$ws.Range('F28').Value = 5
$ws.Range('N28').NumberFormat = "@"
$ws.Range('N28').Value = 'Avela Strumpfhose Belform Liberty Natural  10 - 10.5 6.95 Schweizer Franken'
$ws.Range('O28').NumberFormat = "@"
$ws.Range('O28').Value = '2022-08-18 20:59:13'
$ws.Range('A29').NumberFormat = "@"
$ws.Range('A29').Value = '6075681012'
$ws.Range('B29').NumberFormat = "@"
$ws.Range('B29').Value = 'Avela Strumpfhose Belform Liberty Natural  11 - 12'
$ws.Range('C29').NumberFormat = "@"
$ws.Range('C29').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-belform-liberty-natural-11-12/p/6075681012'
$ws.Range('N29').NumberFormat = "@"
$ws.Range('N29').Value = 'Avela Strumpfhose Belform Liberty Natural  11 - 12 6.95 Schweizer Franken'
$ws.Range('O29').NumberFormat = "@"
$ws.Range('O29').Value = '2022-08-18 20:59:13'
$ws.Range('A30').NumberFormat = "@"
$ws.Range('A30').Value = '6075681009'
$ws.Range('B30').NumberFormat = "@"
$ws.Range('B30').Value = 'Avela Strumpfhose Belform Liberty Natural  8.5 - 9'
$ws.Range('C30').NumberFormat = "@"
$ws.Range('C30').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-belform-liberty-natural-85-9/p/6075681009'
$ws.Range('N30').NumberFormat = "@"
$ws.Range('N30').Value = 'Avela Strumpfhose Belform Liberty Natural  8.5 - 9 6.95 Schweizer Franken'
$ws.Range('O30').NumberFormat = "@"
$ws.Range('O30').Value = '2022-08-18 20:59:13'
$ws.Range('A31').NumberFormat = "@"
$ws.Range('A31').Value = '6075681010'
$ws.Range('B31').NumberFormat = "@"
$ws.Range('B31').Value = 'Avela Strumpfhose Belform Liberty Natural  9.5'
$ws.Range('C31').NumberFormat = "@"
$ws.Range('C31').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-belform-liberty-natural-95/p/6075681010'
$ws.Range('N31').NumberFormat = "@"
$ws.Range('N31').Value = 'Avela Strumpfhose Belform Liberty Natural  9.5 6.95 Schweizer Franken'
$ws.Range('O31').NumberFormat = "@"
$ws.Range('O31').Value = '2022-08-18 20:59:13'
$ws.Range('A32').NumberFormat = "@"
$ws.Range('A32').Value = '6075681015'
$ws.Range('B32').NumberFormat = "@"
$ws.Range('B32').Value = 'Avela Strumpfhose Belform Liberty Noir  10 - 10.5'
$ws.Range('C32').NumberFormat = "@"
$ws.Range('C32').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-belform-liberty-noir-10-105/p/6075681015'
$ws.Range('N32').NumberFormat = "@"
$ws.Range('N32').Value = 'Avela Strumpfhose Belform Liberty Noir  10 - 10.5 6.95 Schweizer Franken'
$ws.Range('O32').NumberFormat = "@"
$ws.Range('O32').Value = '2022-08-18 20:59:13'
$ws.Range('A33').NumberFormat = "@"
$ws.Range('A33').Value = '6075681016'
$ws.Range('B33').NumberFormat = "@"
$ws.Range('B33').Value = 'Avela Strumpfhose Belform Liberty Noir  11 - 12'
$ws.Range('C33').NumberFormat = "@"
$ws.Range('C33').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-belform-liberty-noir-11-12/p/6075681016'
$ws.Range('N33').NumberFormat = "@"
$ws.Range('N33').Value = 'Avela Strumpfhose Belform Liberty Noir  11 - 12 6.95 Schweizer Franken'
$ws.Range('O33').NumberFormat = "@"
$ws.Range('O33').Value = '2022-08-18 20:59:13'
$ws.Range('A34').NumberFormat = "@"
$ws.Range('A34').Value = '6075681013'
$ws.Range('B34').NumberFormat = "@"
$ws.Range('B34').Value = 'Avela Strumpfhose Belform Liberty Noir  8.5 - 9'
$ws.Range('C34').NumberFormat = "@"
$ws.Range('C34').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-belform-liberty-noir-85-9/p/6075681013'
$ws.Range('E34').Value = ""
$ws.Range('F34').Value = 0
$ws.Range('N34').NumberFormat = "@"
$ws.Range('N34').Value = 'Avela Strumpfhose Belform Liberty Noir  8.5 - 9 6.95 Schweizer Franken'
$ws.Range('O34').NumberFormat = "@"
$ws.Range('O34').Value = '2022-08-18 20:59:13'
$ws.Range('A35').NumberFormat = "@"
$ws.Range('A35').Value = '6075681014'
$ws.Range('B35').NumberFormat = "@"
$ws.Range('B35').Value = 'Avela Strumpfhose Belform Liberty Noir  9.5'
$ws.Range('C35').NumberFormat = "@"
$ws.Range('C35').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-belform-liberty-noir-95/p/6075681014'
$ws.Range('E35').Value = 1
$ws.Range('F35').Value = 3
$ws.Range('H35').NumberFormat = "@"
$ws.Range('H35').Value = '6.95'
$ws.Range('N35').NumberFormat = "@"
$ws.Range('N35').Value = 'Avela Strumpfhose Belform Liberty Noir  9.5 6.95 Schweizer Franken'
$ws.Range('O35').NumberFormat = "@"
$ws.Range('O35').Value = '2022-08-18 20:59:13'
$ws.Range('A36').NumberFormat = "@"
$ws.Range('A36').Value = '6075745019'
$ws.Range('B36').NumberFormat = "@"
$ws.Range('B36').Value = 'Avela Strumpfhose Madame 10-10 1/2 nomade'
$ws.Range('C36').NumberFormat = "@"
$ws.Range('C36').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-madame-10-10-12-nomade/p/6075745019'
$ws.Range('N36').NumberFormat = "@"
$ws.Range('N36').Value = 'Avela Strumpfhose Madame 10-10 1/2 nomade 5.95 Schweizer Franken'
$ws.Range('O36').NumberFormat = "@"
$ws.Range('O36').Value = '2022-08-18 20:59:13'
$ws.Range('A37').NumberFormat = "@"
$ws.Range('A37').Value = '6075745020'
$ws.Range('B37').NumberFormat = "@"
$ws.Range('B37').Value = 'Avela Strumpfhose Madame 11-12 nomade'
$ws.Range('C37').NumberFormat = "@"
$ws.Range('C37').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-madame-11-12-nomade/p/6075745020'
$ws.Range('N37').NumberFormat = "@"
$ws.Range('N37').Value = 'Avela Strumpfhose Madame 11-12 nomade 5.95 Schweizer Franken'
$ws.Range('O37').NumberFormat = "@"
$ws.Range('O37').Value = '2022-08-18 20:59:13'
$ws.Range('A38').NumberFormat = "@"
$ws.Range('A38').Value = '6075745018'
$ws.Range('B38').NumberFormat = "@"
$ws.Range('B38').Value = 'Avela Strumpfhose Madame 9- nomade'
$ws.Range('C38').NumberFormat = "@"
$ws.Range('C38').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-madame-9-nomade/p/6075745018'
$ws.Range('N38').NumberFormat = "@"
$ws.Range('N38').Value = 'Avela Strumpfhose Madame 9- nomade 5.95 Schweizer Franken'
$ws.Range('O38').NumberFormat = "@"
$ws.Range('O38').Value = '2022-08-18 20:59:13'
$ws.Range('A39').NumberFormat = "@"
$ws.Range('A39').Value = '6075745008'
$ws.Range('B39').NumberFormat = "@"
$ws.Range('B39').Value = 'Avela Strumpfhose Madame Hasel  11 - 12'
$ws.Range('C39').NumberFormat = "@"
$ws.Range('C39').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-madame-hasel-11-12/p/6075745008'
$ws.Range('N39').NumberFormat = "@"
$ws.Range('N39').Value = 'Avela Strumpfhose Madame Hasel  11 - 12 5.95 Schweizer Franken'
$ws.Range('O39').NumberFormat = "@"
$ws.Range('O39').Value = '2022-08-18 20:59:13'
$ws.Range('A40').NumberFormat = "@"
$ws.Range('A40').Value = '6075745005'
$ws.Range('B40').NumberFormat = "@"
$ws.Range('B40').Value = 'Avela Strumpfhose Madame Hasel  8.5 - 9'
$ws.Range('C40').NumberFormat = "@"
$ws.Range('C40').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-madame-hasel-85-9/p/6075745005'
$ws.Range('N40').NumberFormat = "@"
$ws.Range('N40').Value = 'Avela Strumpfhose Madame Hasel  8.5 - 9 5.95 Schweizer Franken'
$ws.Range('O40').NumberFormat = "@"
$ws.Range('O40').Value = '2022-08-18 20:59:13'
$ws.Range('A41').NumberFormat = "@"
$ws.Range('A41').Value = '6075745017'
$ws.Range('B41').NumberFormat = "@"
$ws.Range('B41').Value = 'Avela Strumpfhose Madame Nomade  8.5 - 9'
$ws.Range('C41').NumberFormat = "@"
$ws.Range('C41').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-madame-nomade-85-9/p/6075745017'
$ws.Range('H41').NumberFormat = "@"
$ws.Range('H41').Value = '5.95'
$ws.Range('N41').NumberFormat = "@"
$ws.Range('N41').Value = 'Avela Strumpfhose Madame Nomade  8.5 - 9 5.95 Schweizer Franken'
$ws.Range('O41').NumberFormat = "@"
$ws.Range('O41').Value = '2022-08-18 20:59:13'
$ws.Range('A42').NumberFormat = "@"
$ws.Range('A42').Value = '6076125011'
$ws.Range('B42').NumberFormat = "@"
$ws.Range('B42').Value = 'Avela Strumpfhose Top Silhouette Natural  10 - 10.5'
$ws.Range('C42').NumberFormat = "@"
$ws.Range('C42').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-top-silhouette-natural-10-105/p/6076125011'
$ws.Range('N42').NumberFormat = "@"
$ws.Range('N42').Value = 'Avela Strumpfhose Top Silhouette Natural  10 - 10.5 9.95 Schweizer Franken'
$ws.Range('O42').NumberFormat = "@"
$ws.Range('O42').Value = '2022-08-18 20:59:13'
$ws.Range('A43').NumberFormat = "@"
$ws.Range('A43').Value = '6076125012'
$ws.Range('B43').NumberFormat = "@"
$ws.Range('B43').Value = 'Avela Strumpfhose Top Silhouette Natural  11 - 12'
$ws.Range('C43').NumberFormat = "@"
$ws.Range('C43').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-top-silhouette-natural-11-12/p/6076125012'
$ws.Range('N43').NumberFormat = "@"
$ws.Range('N43').Value = 'Avela Strumpfhose Top Silhouette Natural  11 - 12 9.95 Schweizer Franken'
$ws.Range('O43').NumberFormat = "@"
$ws.Range('O43').Value = '2022-08-18 20:59:13'
$ws.Range('A44').NumberFormat = "@"
$ws.Range('A44').Value = '6076125009'
$ws.Range('B44').NumberFormat = "@"
$ws.Range('B44').Value = 'Avela Strumpfhose Top Silhouette Natural  8.5 - 9'
$ws.Range('C44').NumberFormat = "@"
$ws.Range('C44').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-top-silhouette-natural-85-9/p/6076125009'
$ws.Range('E44').Value = ""
$ws.Range('F44').Value = 0
$ws.Range('N44').NumberFormat = "@"
$ws.Range('N44').Value = 'Avela Strumpfhose Top Silhouette Natural  8.5 - 9 9.95 Schweizer Franken'
$ws.Range('O44').NumberFormat = "@"
$ws.Range('O44').Value = '2022-08-18 20:59:13'
$ws.Range('A45').NumberFormat = "@"
$ws.Range('A45').Value = '6076125010'
$ws.Range('B45').NumberFormat = "@"
$ws.Range('B45').Value = 'Avela Strumpfhose Top Silhouette Natural  9.5'
$ws.Range('C45').NumberFormat = "@"
$ws.Range('C45').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-top-silhouette-natural-95/p/6076125010'
$ws.Range('N45').NumberFormat = "@"
$ws.Range('N45').Value = 'Avela Strumpfhose Top Silhouette Natural  9.5 9.95 Schweizer Franken'
$ws.Range('O45').NumberFormat = "@"
$ws.Range('O45').Value = '2022-08-18 20:59:13'
$ws.Range('A46').NumberFormat = "@"
$ws.Range('A46').Value = '6076125015'
$ws.Range('B46').NumberFormat = "@"
$ws.Range('B46').Value = 'Avela Strumpfhose Top Silhouette Noir  10 - 10.5'
$ws.Range('C46').NumberFormat = "@"
$ws.Range('C46').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-top-silhouette-noir-10-105/p/6076125015'
$ws.Range('E46').Value = 1
$ws.Range('F46').Value = 5
$ws.Range('N46').NumberFormat = "@"
$ws.Range('N46').Value = 'Avela Strumpfhose Top Silhouette Noir  10 - 10.5 9.95 Schweizer Franken'
$ws.Range('O46').NumberFormat = "@"
$ws.Range('O46').Value = '2022-08-18 20:59:13'
$ws.Range('A47').NumberFormat = "@"
$ws.Range('A47').Value = '6076125016'
$ws.Range('B47').NumberFormat = "@"
$ws.Range('B47').Value = 'Avela Strumpfhose Top Silhouette Noir  11 - 12'
$ws.Range('C47').NumberFormat = "@"
$ws.Range('C47').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-top-silhouette-noir-11-12/p/6076125016'
$ws.Range('N47').NumberFormat = "@"
$ws.Range('N47').Value = 'Avela Strumpfhose Top Silhouette Noir  11 - 12 9.95 Schweizer Franken'
$ws.Range('O47').NumberFormat = "@"
$ws.Range('O47').Value = '2022-08-18 20:59:13'
$ws.Range('A48').NumberFormat = "@"
$ws.Range('A48').Value = '6076125013'
$ws.Range('B48').NumberFormat = "@"
$ws.Range('B48').Value = 'Avela Strumpfhose Top Silhouette Noir  8.5 - 9'
$ws.Range('C48').NumberFormat = "@"
$ws.Range('C48').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-top-silhouette-noir-85-9/p/6076125013'
$ws.Range('E48').Value = ""
$ws.Range('F48').Value = 0
$ws.Range('N48').NumberFormat = "@"
$ws.Range('N48').Value = 'Avela Strumpfhose Top Silhouette Noir  8.5 - 9 9.95 Schweizer Franken'
$ws.Range('O48').NumberFormat = "@"
$ws.Range('O48').Value = '2022-08-18 20:59:13'
$ws.Range('A49').NumberFormat = "@"
$ws.Range('A49').Value = '6076125014'
$ws.Range('B49').NumberFormat = "@"
$ws.Range('B49').Value = 'Avela Strumpfhose Top Silhouette Noir  9.5'
$ws.Range('C49').NumberFormat = "@"
$ws.Range('C49').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-top-silhouette-noir-95/p/6076125014'
$ws.Range('E49').Value = 1
$ws.Range('F49').Value = 5
$ws.Range('G49').NumberFormat = "@"
$ws.Range('G49').Value = 'Avela'
$ws.Range('H49').NumberFormat = "@"
$ws.Range('H49').Value = '9.95'
$ws.Range('M49').NumberFormat = "@"
$ws.Range('M49').Value = '[''haushalt-tier'', ''bekleidung'', ''socken-unterwaesche'', ''struempfe'']'
$ws.Range('N49').NumberFormat = "@"
$ws.Range('N49').Value = 'Avela Strumpfhose Top Silhouette Noir  9.5 9.95 Schweizer Franken'
$ws.Range('O49').NumberFormat = "@"
$ws.Range('O49').Value = '2022-08-18 20:59:13'
$ws.Range('A50').NumberFormat = "@"
$ws.Range('A50').Value = '5799901002'
$ws.Range('B50').NumberFormat = "@"
$ws.Range('B50').Value = 'Magic Matic schwarz ecorepel'
$ws.Range('C50').NumberFormat = "@"
$ws.Range('C50').Value = '/de/haushalt-tier/bekleidung/taschen-accessoires/schirme/magic-matic-schwarz-ecorepel/p/5799901002'
$ws.Range('E50').Value = 3
$ws.Range('N50').NumberFormat = "@"
$ws.Range('N50').Value = 'Magic Matic schwarz ecorepel 19.95 Schweizer Franken'
$ws.Range('O50').NumberFormat = "@"
$ws.Range('O50').Value = '2022-08-18 20:59:13'
$ws.Range('A51').NumberFormat = "@"
$ws.Range('A51').Value = '5799902001'
$ws.Range('B51').NumberFormat = "@"
$ws.Range('B51').Value = 'Mini Matic schwarz ecorepel'
$ws.Range('C51').NumberFormat = "@"
$ws.Range('C51').Value = '/de/haushalt-tier/bekleidung/taschen-accessoires/schirme/mini-matic-schwarz-ecorepel/p/5799902001'
$ws.Range('E51').Value = 1
$ws.Range('F51').Value = 1
$ws.Range('H51').NumberFormat = "@"
$ws.Range('H51').Value = '19.95'
$ws.Range('M51').NumberFormat = "@"
$ws.Range('M51').Value = '[''haushalt-tier'', ''bekleidung'', ''taschen-accessoires'', ''schirme'']'
$ws.Range('N51').NumberFormat = "@"
$ws.Range('N51').Value = 'Mini Matic schwarz ecorepel 19.95 Schweizer Franken'
$ws.Range('O51').NumberFormat = "@"
$ws.Range('O51').Value = '2022-08-18 20:59:13'
$ws.Range('A52').NumberFormat = "@"
$ws.Range('A52').Value = '3875554004'
$ws.Range('B52').NumberFormat = "@"
$ws.Range('B52').Value = 'Naturaline Damen Bustier Schwarz S'
$ws.Range('C52').NumberFormat = "@"
$ws.Range('C52').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/unterwaesche/naturaline-damen-bustier-schwarz-s/p/3875554004'
$ws.Range('N52').NumberFormat = "@"
$ws.Range('N52').Value = 'Naturaline Damen Bustier Schwarz S 14.95 Schweizer Franken'
$ws.Range('O52').NumberFormat = "@"
$ws.Range('O52').Value = '2022-08-18 20:59:13'
$ws.Range('A53').NumberFormat = "@"
$ws.Range('A53').Value = '3875554007'
$ws.Range('B53').NumberFormat = "@"
$ws.Range('B53').Value = 'Naturaline Damen Bustier Weiss S'
$ws.Range('C53').NumberFormat = "@"
$ws.Range('C53').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/unterwaesche/naturaline-damen-bustier-weiss-s/p/3875554007'
$ws.Range('E53').Value = ""
$ws.Range('F53').Value = 0
$ws.Range('H53').NumberFormat = "@"
$ws.Range('H53').Value = '14.95'
$ws.Range('N53').NumberFormat = "@"
$ws.Range('N53').Value = 'Naturaline Damen Bustier Weiss S 14.95 Schweizer Franken'
$ws.Range('O53').NumberFormat = "@"
$ws.Range('O53').Value = '2022-08-18 20:59:13'
$ws.Range('A54').NumberFormat = "@"
$ws.Range('A54').Value = '3305779007'
$ws.Range('B54').NumberFormat = "@"
$ws.Range('B54').Value = 'Naturaline Damen Panty S weiss'
$ws.Range('C54').NumberFormat = "@"
$ws.Range('C54').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/unterwaesche/naturaline-damen-panty-s-weiss/p/3305779007'
$ws.Range('E54').Value = 1
$ws.Range('F54').Value = 4
$ws.Range('N54').NumberFormat = "@"
$ws.Range('N54').Value = 'Naturaline Damen Panty S weiss 9.95 Schweizer Franken'
$ws.Range('O54').NumberFormat = "@"
$ws.Range('O54').Value = '2022-08-18 20:59:13'
$ws.Range('A55').NumberFormat = "@"
$ws.Range('A55').Value = '4322745004'
$ws.Range('B55').NumberFormat = "@"
$ws.Range('B55').Value = 'Naturaline Damen String schwarz L'
$ws.Range('C55').NumberFormat = "@"
$ws.Range('C55').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/unterwaesche/naturaline-damen-string-schwarz-l/p/4322745004'
$ws.Range('N55').NumberFormat = "@"
$ws.Range('N55').Value = 'Naturaline Damen String schwarz L 9.95 Schweizer Franken'
$ws.Range('O55').NumberFormat = "@"
$ws.Range('O55').Value = '2022-08-18 20:59:13'
$ws.Range('A56').NumberFormat = "@"
$ws.Range('A56').Value = '4322745002'
$ws.Range('B56').NumberFormat = "@"
$ws.Range('B56').Value = 'Naturaline Damen String schwarz S'
$ws.Range('C56').NumberFormat = "@"
$ws.Range('C56').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/unterwaesche/naturaline-damen-string-schwarz-s/p/4322745002'
$ws.Range('N56').NumberFormat = "@"
$ws.Range('N56').Value = 'Naturaline Damen String schwarz S 9.95 Schweizer Franken'
$ws.Range('O56').NumberFormat = "@"
$ws.Range('O56').Value = '2022-08-18 20:59:13'
$ws.Range('A57').NumberFormat = "@"
$ws.Range('A57').Value = '4322745007'
$ws.Range('B57').NumberFormat = "@"
$ws.Range('B57').Value = 'Naturaline Damen String weiss M'
$ws.Range('C57').NumberFormat = "@"
$ws.Range('C57').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/unterwaesche/naturaline-damen-string-weiss-m/p/4322745007'
$ws.Range('H57').NumberFormat = "@"
$ws.Range('H57').Value = '9.95'
$ws.Range('N57').NumberFormat = "@"
$ws.Range('N57').Value = 'Naturaline Damen String weiss M 9.95 Schweizer Franken'
$ws.Range('O57').NumberFormat = "@"
$ws.Range('O57').Value = '2022-08-18 20:59:13'
$ws.Range('A58').NumberFormat = "@"
$ws.Range('A58').Value = '3404677005'
$ws.Range('B58').NumberFormat = "@"
$ws.Range('B58').Value = 'Naturaline Herren Slip schwarz S'
$ws.Range('C58').NumberFormat = "@"
$ws.Range('C58').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/unterwaesche/naturaline-herren-slip-schwarz-s/p/3404677005'
$ws.Range('H58').NumberFormat = "@"
$ws.Range('H58').Value = '14.95'
$ws.Range('N58').NumberFormat = "@"
$ws.Range('N58').Value = 'Naturaline Herren Slip schwarz S 14.95 Schweizer Franken'
$ws.Range('O58').NumberFormat = "@"
$ws.Range('O58').Value = '2022-08-18 20:59:13'
$ws.Range('A59').NumberFormat = "@"
$ws.Range('A59').Value = '3305289015'
$ws.Range('B59').NumberFormat = "@"
$ws.Range('B59').Value = 'Naturaline Herren Slip weiss L'
$ws.Range('C59').NumberFormat = "@"
$ws.Range('C59').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/unterwaesche/naturaline-herren-slip-weiss-l/p/3305289015'
$ws.Range('N59').NumberFormat = "@"
$ws.Range('N59').Value = 'Naturaline Herren Slip weiss L 9.95 Schweizer Franken'
$ws.Range('O59').NumberFormat = "@"
$ws.Range('O59').Value = '2022-08-18 20:59:13'
$ws.Range('A60').NumberFormat = "@"
$ws.Range('A60').Value = '3305289014'
$ws.Range('B60').NumberFormat = "@"
$ws.Range('B60').Value = 'Naturaline Herren Slip weiss M'
$ws.Range('C60').NumberFormat = "@"
$ws.Range('C60').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/unterwaesche/naturaline-herren-slip-weiss-m/p/3305289014'
$ws.Range('N60').NumberFormat = "@"
$ws.Range('N60').Value = 'Naturaline Herren Slip weiss M 9.95 Schweizer Franken'
$ws.Range('O60').NumberFormat = "@"
$ws.Range('O60').Value = '2022-08-18 20:59:13'
$ws.Range('A61').NumberFormat = "@"
$ws.Range('A61').Value = '3305289013'
$ws.Range('B61').NumberFormat = "@"
$ws.Range('B61').Value = 'Naturaline Herren Slip weiss S'
$ws.Range('C61').NumberFormat = "@"
$ws.Range('C61').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/unterwaesche/naturaline-herren-slip-weiss-s/p/3305289013'
$ws.Range('E61').Value = ""
$ws.Range('F61').Value = 0
$ws.Range('N61').NumberFormat = "@"
$ws.Range('N61').Value = 'Naturaline Herren Slip weiss S 9.95 Schweizer Franken'
$ws.Range('O61').NumberFormat = "@"
$ws.Range('O61').Value = '2022-08-18 20:59:13'
$ws.Range('A62').NumberFormat = "@"
$ws.Range('A62').Value = '3305289016'
$ws.Range('B62').NumberFormat = "@"
$ws.Range('B62').Value = 'Naturaline Herren Slip weiss XL'
$ws.Range('C62').NumberFormat = "@"
$ws.Range('C62').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/unterwaesche/naturaline-herren-slip-weiss-xl/p/3305289016'
$ws.Range('E62').Value = 1
$ws.Range('F62').Value = 2
$ws.Range('M62').NumberFormat = "@"
$ws.Range('M62').Value = '[''haushalt-tier'', ''bekleidung'', ''socken-unterwaesche'', ''unterwaesche'']'
$ws.Range('N62').NumberFormat = "@"
$ws.Range('N62').Value = 'Naturaline Herren Slip weiss XL 9.95 Schweizer Franken'
$ws.Range('O62').NumberFormat = "@"
$ws.Range('O62').Value = '2022-08-18 20:59:13'
$ws.Range('A63').NumberFormat = "@"
$ws.Range('A63').Value = '6365980001'
$ws.Range('B63').NumberFormat = "@"
$ws.Range('B63').Value = 'Naturaline Herren Socken Glatt Duo Dunkelgrau 40 - 42'
$ws.Range('C63').NumberFormat = "@"
$ws.Range('C63').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/socken/naturaline-herren-socken-glatt-duo-dunkelgrau-40-42/p/6365980001'
$ws.Range('N63').NumberFormat = "@"
$ws.Range('N63').Value = 'Naturaline Herren Socken Glatt Duo Dunkelgrau 40 - 42 9.95 Schweizer Franken'
$ws.Range('O63').NumberFormat = "@"
$ws.Range('O63').Value = '2022-08-18 20:59:13'
$ws.Range('A64').NumberFormat = "@"
$ws.Range('A64').Value = '6365980011'
$ws.Range('B64').NumberFormat = "@"
$ws.Range('B64').Value = 'Naturaline Herren Socken Glatt Duo Weiss 40 - 42'
$ws.Range('C64').NumberFormat = "@"
$ws.Range('C64').Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/socken/naturaline-herren-socken-glatt-duo-weiss-40-42/p/6365980011'
$ws.Range('H64').NumberFormat = "@"
$ws.Range('H64').Value = '9.95'
$ws.Range('M64').NumberFormat = "@"
$ws.Range('M64').Value = '[''haushalt-tier'', ''bekleidung'', ''socken-unterwaesche'', ''socken'']'
$ws.Range('N64').NumberFormat = "@"
$ws.Range('N64').Value = 'Naturaline Herren Socken Glatt Duo Weiss 40 - 42 9.95 Schweizer Franken'
$ws.Range('O64').NumberFormat = "@"
$ws.Range('O64').Value = '2022-08-18 20:59:13'
$ws.Range('A65').NumberFormat = "@"
$ws.Range('A65').Value = '6031467016'
$ws.Range('B65').NumberFormat = "@"
$ws.Range('B65').Value = 'Naturaline Herren T-Shirt Kurzarm schwarz S'
$ws.Range('C65').NumberFormat = "@"
$ws.Range('C65').Value = '/de/haushalt-tier/bekleidung/shirts-pullover/herren-shirt/naturaline-herren-t-shirt-kurzarm-schwarz-s/p/6031467016'
$ws.Range('N65').NumberFormat = "@"
$ws.Range('N65').Value = 'Naturaline Herren T-Shirt Kurzarm schwarz S 24.95 Schweizer Franken'
$ws.Range('O65').NumberFormat = "@"
$ws.Range('O65').Value = '2022-08-18 20:59:13'
$ws.Range('A66').NumberFormat = "@"
$ws.Range('A66').Value = '6031467009'
$ws.Range('B66').NumberFormat = "@"
$ws.Range('B66').Value = 'Naturaline Herren T-Shirt Kurzarm weissXL'
$ws.Range('C66').NumberFormat = "@"
$ws.Range('C66').Value = '/de/haushalt-tier/bekleidung/shirts-pullover/herren-shirt/naturaline-herren-t-shirt-kurzarm-weissxl/p/6031467009'
$ws.Range('N66').NumberFormat = "@"
$ws.Range('N66').Value = 'Naturaline Herren T-Shirt Kurzarm weissXL - Online kein Bestand 24.95 Schweizer Franken'
$ws.Range('O66').NumberFormat = "@"
$ws.Range('O66').Value = '2022-08-18 20:59:13'
$ws.Range('A67').NumberFormat = "@"
$ws.Range('A67').Value = '6031467010'
$ws.Range('B67').NumberFormat = "@"
$ws.Range('B67').Value = 'Naturaline Herren T-Shirt Kurzarm weiss XXL'
$ws.Range('C67').NumberFormat = "@"
$ws.Range('C67').Value = '/de/haushalt-tier/bekleidung/shirts-pullover/herren-shirt/naturaline-herren-t-shirt-kurzarm-weiss-xxl/p/6031467010'
$ws.Range('H67').NumberFormat = "@"
$ws.Range('H67').Value = '24.95'
$ws.Range('M67').NumberFormat = "@"
$ws.Range('M67').Value = '[''haushalt-tier'', ''bekleidung'', ''shirts-pullover'', ''herren-shirt'']'
$ws.Range('N67').NumberFormat = "@"
$ws.Range('N67').Value = 'Naturaline Herren T-Shirt Kurzarm weiss XXL 24.95 Schweizer Franken'
$ws.Range('O67').NumberFormat = "@"
$ws.Range('O67').Value = '2022-08-18 20:59:13'
$ws.Range('A68').NumberFormat = "@"
$ws.Range('A68').Value = '3890690001'
$ws.Range('B68').NumberFormat = "@"
$ws.Range('B68').Value = 'Portemonnaie schwarz quer klein'
$ws.Range('C68').NumberFormat = "@"
$ws.Range('C68').Value = '/de/haushalt-tier/bekleidung/taschen-accessoires/portemonnaie/portemonnaie-schwarz-quer-klein/p/3890690001'
$ws.Range('D68').Value = ""
$ws.Range('G68').NumberFormat = "@"
$ws.Range('G68').Value = 'Coop'
$ws.Range('H68').NumberFormat = "@"
$ws.Range('H68').Value = '34.95'
$ws.Range('I68').Value = ""
$ws.Range('J68').Value = ""
$ws.Range('K68').Value = ""
$ws.Range('L68').Value = ""
$ws.Range('M68').NumberFormat = "@"
$ws.Range('M68').Value = '[''haushalt-tier'', ''bekleidung'', ''taschen-accessoires'', ''portemonnaie'']'
$ws.Range('N68').NumberFormat = "@"
$ws.Range('N68').Value = 'Portemonnaie schwarz quer klein 34.95 Schweizer Franken'
$ws.Range('O68').NumberFormat = "@"
$ws.Range('O68').Value = '2022-08-18 20:59:13'
$ws.Range('O69').NumberFormat = "@"
$ws.Range('O69').Value = '2022-08-18 20:59:13'
$ws.Range('O70').NumberFormat = "@"
$ws.Range('O70').Value = '2022-08-18 20:59:13'
$ws.Range('O71').NumberFormat = "@"
$ws.Range('O71').Value = '2022-08-18 20:59:13'
$ws.Range('O72').NumberFormat = "@"
$ws.Range('O72').Value = '2022-08-18 20:59:13'
$ws.Range('O73').NumberFormat = "@"
$ws.Range('O73').Value = '2022-08-18 20:59:13'
